$wb = $excel.ActiveWorkbook

# ----- "Immediate Checklist" sheet (sheet3): refresh the to-do items -----
$todo = $wb.Worksheets.Item("Immediate Checklist")

$todo.Range("A3").Value = "Play a different sound in the shop when ammo cannot be purchased either due to reaching the max ammo or not having enough money."
$todo.Range("A3").WrapText = $true
$todo.Rows.Item(3).RowHeight = 30

$todo.Range("A4").Value = "Ammo starts on 0/0 for some reason in builds."
$todo.Range("A4").WrapText = $true

$todo.Range("A5").Value = "Lock the enemy health bars to only rotate on y-axis."
$todo.Range("A5").WrapText = $true

$todo.Range("A6").Value = "Dealing 100 damange with sniper rifle to a 20 HP enemy reduces wave HP counter by 100 instead of 20, resulting in negative values."
$todo.Range("A6").WrapText = $true
$todo.Rows.Item(6).RowHeight = 30

# New blank row 7 (keeps the same wrap-text formatting as the rest of the column)
$todo.Range("A7").WrapText = $true

# ----- "Features" sheet (sheet1): mark A12 ("Add aim down sights animation") as Neutral -----
$features = $wb.Worksheets.Item("Features")
$features.Range("A12").Style = "Neutral"
$features.Range("A12").WrapText = $true

# ----- "Bugs and Quirks" sheet (sheet2): no content changes, just selection -----
$bugs = $wb.Worksheets.Item("Bugs and Quirks")

# ----- Restore / update the per-sheet selections -----
$features.Range("A12").Select()
$bugs.Range("A3").Select()
$todo.Range("A7").Select()
